$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row at the TOP of the table (new row 4) for
#    "ALBENDAZOLE 400MG/5ML SUSP. 30ML". This pushes all existing data rows,
#    the totals row and the footer row down by one.
#    (Insert first, THEN copy the sibling row's formatting into the newly
#    created blank row - doing it the other way round loses the clipboard
#    contents across the Insert call.)
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()
$ws.Range("A5:N5").Copy()
$ws.Range("A4:N4").PasteSpecial(-4122)   # xlPasteFormats - reuse existing styles
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(5).RowHeight()

$ws.Range("B4:G4").Merge()
$ws.Range("H4:K4").Merge()
$ws.Range("L4:M4").Merge()

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "ALBENDAZOLE 400MG/5ML SUSP. 30ML"
$ws.Range("H4").Value = "4:0"
$ws.Range("L4").Value = 48
$ws.Range("N4").Value = "2:0"

# Renumber the sequence column for the rows that used to be 4-8 and are now
# 5-9 (their content/styles already shifted down intact via the row insert).
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6

# ---------------------------------------------------------------------------
# 2) Insert a new data row at the BOTTOM of the table (after the now-shifted
#    last data row, row 9) for "VISCERALGINE 10MG/5ML SYRUP 120 ML".
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).Insert()
$ws.Range("A9:N9").Copy()
$ws.Range("A10:N10").PasteSpecial(-4122)   # xlPasteFormats - reuse existing styles
$ws.Rows.Item(10).RowHeight = $ws.Rows.Item(9).RowHeight()

$ws.Range("B10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "VISCERALGINE 10MG/5ML SYRUP 120 ML"
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").Value = 35
$ws.Range("N10").Value = "1:0"

# ---------------------------------------------------------------------------
# 3) Update the grand total (now on row 11, was row 9) to reflect the two
#    newly added rows: 146 + 48 + 35 = 229.
# ---------------------------------------------------------------------------
$ws.Range("K11").Value = 229
